$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new price value is numeric-looking: pre-format as Text
# so Excel stores the literal string instead of coercing to a Number.
$textCells = @("D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D14", "D15", "D16", "D17", "D19", "D20", "D21", "D22", "D24", "D25", "D28", "D30", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range('D2').Value = '30.700.86'
$ws.Range('E2').Value = '  +0.76%  '
$ws.Range('D3').Value = '2.123.91'
$ws.Range('E3').Value = '  +1.03%  '
$ws.Range('E4').Value = '  +1.14%  '
$ws.Range('D5').Value = '339.00'
$ws.Range('D6').Value = '1.012'
$ws.Range('E6').Value = '  +1.01%  '
$ws.Range('D7').Value = '0.5278'
$ws.Range('E7').Value = '  +1.00%  '
$ws.Range('D8').Value = '0.4573'
$ws.Range('E8').Value = '  +1.72%  '
$ws.Range('D9').Value = '54.54'
$ws.Range('E9').Value = '  +1.27%  '
$ws.Range('D10').Value = '0.09132'
$ws.Range('E10').Value = '  +2.34%  '
$ws.Range('D11').Value = '1.177'
$ws.Range('E11').Value = '  +1.83%  '
$ws.Range('D12').Value = '24.60'
$ws.Range('E12').Value = '  +1.09%  '
$ws.Range('D13').Value = '2.122.59'
$ws.Range('E13').Value = '  +1.62%  '
$ws.Range('D14').Value = '6.873'
$ws.Range('E14').Value = '  +2.07%  '
$ws.Range('D15').Value = '8.127'
$ws.Range('E15').Value = '  +5.19%  '
$ws.Range('D16').Value = '0.00001178'
$ws.Range('E16').Value = '  +4.70%  '
$ws.Range('D17').Value = '97.39'
$ws.Range('E17').Value = '  +1.04%  '
$ws.Range('D19').Value = '0.06709'
$ws.Range('E19').Value = '  +1.40%  '
$ws.Range('D20').Value = '19.62'
$ws.Range('E20').Value = '  +2.10%  '
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').Value = '6.519'
$ws.Range('E21').Value = '  +3.67%  '
$ws.Range('B22').Value = 'Dai'
$ws.Range('C22').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D22').Value = '1.012'
$ws.Range('E22').Value = '  +0.96%  '
$ws.Range('D23').Value = '30.773.23'
$ws.Range('E23').Value = '  +0.80%  '
$ws.Range('D24').Value = '13.03'
$ws.Range('E24').Value = '  +5.75%  '
$ws.Range('D25').Value = '2.361'
$ws.Range('E25').Value = '  +1.75%  '
$ws.Range('D26').Value = '2.376.81'
$ws.Range('E26').Value = '  +1.72%  '
$ws.Range('E27').Value = '  +1.03%  '
$ws.Range('D28').Value = '165.53'
$ws.Range('E28').Value = '  +1.06%  '
$ws.Range('E29').Value = '  -0.72%  '
$ws.Range('D30').Value = '135.84'
$ws.Range('E30').Value = '  +2.45%  '
$ws.Range('E31').Value = '  +1.06%  '
$ws.Range('D32').Value = '0.1081'
$ws.Range('E32').Value = '  +0.76%  '
$ws.Range('D33').Value = '1.658'
$ws.Range('E33').Value = '  -0.61%  '
$ws.Range('D34').Value = '6.432'
$ws.Range('E34').Value = '  +4.67%  '
$ws.Range('D35').Value = '3.947'
$ws.Range('E35').Value = '  +0.32%  '
$ws.Range('D36').Value = '10.62'
$ws.Range('E36').Value = '  +1.78%  '
$ws.Range('D37').Value = '5.964'
$ws.Range('E37').Value = '  +8.84%  '
$ws.Range('D38').Value = '0.02689'
$ws.Range('E38').Value = '  +4.66%  '
$ws.Range('D39').Value = '0.06909'
$ws.Range('E39').Value = '  +1.99%  '
$ws.Range('D40').Value = '0.2336'
$ws.Range('E40').Value = '  +3.06%  '
$ws.Range('D41').Value = '12.68'
$ws.Range('E41').Value = '  -0.98%  '
$ws.Range('D42').Value = '0.6942'
$ws.Range('E42').Value = '  +0.29%  '
$ws.Range('E43').Value = '  +0.79%  '
$ws.Range('D44').Value = '15.19'
$ws.Range('E44').Value = '  +8.34%  '
$ws.Range('D45').Value = '0.6505'
$ws.Range('E45').Value = '  +2.14%  '
$ws.Range('D46').Value = '2.319'
$ws.Range('E46').Value = '  +0.71%  '
$ws.Range('D47').Value = '0.00000000371'
$ws.Range('E47').Value = '  +16.24%  '
$ws.Range('D48').Value = '3.703'
$ws.Range('E48').Value = '  +1.79%  '
$ws.Range('D49').Value = '1.262'
$ws.Range('E49').Value = '  +1.47%  '
$ws.Range('D50').Value = '84.03'
$ws.Range('E50').Value = '  +1.24%  '
$ws.Range('D51').Value = '0.07315'
$ws.Range('E51').Value = '  +3.83%  '
